# Regenerate save_data to use K (strikeouts) instead of Strike# (pitch count proxy),
# rewriting the computed "K" column (G) values for each start.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new K value (col G), per regenerated std/mean/s_vals calculation.
$kValues = @{
    2  = 2
    3  = 0
    4  = 0
    5  = 0
    6  = 1
    7  = 3
    8  = 5
    9  = 5
    10 = 1
    11 = 4
    12 = 3
    13 = 1
    14 = 3
    15 = 6
    16 = 4
    17 = 1
    18 = 1
    19 = 3
    20 = 3
    21 = 5
    22 = 0
    23 = 4
    24 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
